$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column headers were renamed to match the PCB "mid point" terminology.
$ws.Range("B1").Value = "MidX"
$ws.Range("C1").Value = "MidY"

# Corrected rotation for U1 (was 0, should be 270).
$ws.Range("D6").Value = 270

# Clear the stray numeric-format styling that was left on the Designator/Layer columns.
$ws.Range("A2:A6").Style = "Normal"
$ws.Range("E2:E6").Style = "Normal"

# Restore the cursor/selection and window position as last left by the author.
$ws.Range("D7").Select()
$win = $excel.ActiveWindow
$win.Left = 6960
$win.Top = 6024
